$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.38094952709161589
$ws.Range("B1").Value = 0.37984142217757721
$ws.Range("A2").Value = -0.27435413431515698
$ws.Range("B2").Value = 0.27124761117078044
$ws.Range("A3").Value = -0.16829455490639944
$ws.Range("B3").Value = 0.16736047688051769
$ws.Range("A4").Value = -0.15536047698250144
$ws.Range("B4").Value = 0.15452405893862498
$ws.Range("A5").Value = -0.14852405930034518
$ws.Range("B5").Value = 0.14685179533850246
$ws.Range("A6").Value = -0.1094531974645121
$ws.Range("B6").Value = 0.10926772975338395
$ws.Range("A7").Value = -0.089267730200988993
$ws.Range("B7").Value = 0.088799924077182624
$ws.Range("A8").Value = -0.06879992452986361
$ws.Range("B8").Value = 0.068395307561572416
$ws.Range("A9").Value = -0.062395307947272549
$ws.Range("B9").Value = 0.062050723510822969
$ws.Range("A10").Value = -0.056050723902167476
$ws.Range("B10").Value = 0.056001570621781127
$ws.Range("A11").Value = -0.051501571006102864
$ws.Range("B11").Value = 0.051418354659137577
$ws.Range("A12").Value = -0.045418355052796233
$ws.Range("B12").Value = 0.045160202343113909
$ws.Range("A13").Value = -0.039160202742775319
$ws.Range("B13").Value = 0.039089905466278552
$ws.Range("A14").Value = -0.027089905899178035
$ws.Range("B14").Value = 0.027055919735624379
$ws.Range("A15").Value = -0.021055920138508988
$ws.Range("B15").Value = 0.02102916425946777
$ws.Range("A16").Value = -0.015029164663616479
$ws.Range("B16").Value = 0.015004687577876208
$ws.Range("A17").Value = -0.0090046879837490934
$ws.Range("B17").Value = 0.008999999578027662
$ws.Range("A18").Value = -0.095487288224859412
$ws.Range("B18").Value = 0.095363788258133297
$ws.Range("A19").Value = -0.086363788616273585
$ws.Range("B19").Value = 0.085401433816318928
$ws.Range("A20").Value = -0.076401434185432215
$ws.Range("B20").Value = 0.076188991606175094
$ws.Range("A21").Value = -0.0090044246681810769
$ws.Range("B21").Value = 0.0089999996280858419
$ws.Range("A22").Value = -0.093953056632537724
$ws.Range("B22").Value = 0.093638198281123408
$ws.Range("A23").Value = -0.084638198652578822
$ws.Range("B23").Value = 0.084127593368980236
$ws.Range("A24").Value = -0.042127593918869444
$ws.Range("B24").Value = 0.041999999447164704
$ws.Range("A25").Value = -0.031398598246553888
$ws.Range("B25").Value = 0.031367741876493938
$ws.Range("A26").Value = -0.009725443730758343
$ws.Range("B26").Value = 0.009708448786657442
$ws.Range("A27").Value = -0.0037084491598928793
$ws.Range("B27").Value = 0.003652433664995236
$ws.Range("A28").Value = 0.002347565961437148
$ws.Range("B28").Value = -0.0023836468866056393
$ws.Range("A29").Value = 0.014383646481828549
$ws.Range("B29").Value = -0.01440163662062055
$ws.Range("A30").Value = 0.034401636174621331
$ws.Range("B30").Value = -0.034591973157241274
$ws.Range("A31").Value = 0.049591972740335777
$ws.Range("B31").Value = -0.049691030944529757
$ws.Range("A32").Value = -0.0060008110496596601
$ws.Range("B32").Value = 0.0059999996313200299
